$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.927.75'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '2.361.92'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('E5').Value = '  +6.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.27'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '76.20'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +6.42%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.625'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +25.38%  '
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.27'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('E12').Value = '  +19.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +18.88%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '2.715.59'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '16.88'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +4.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.918'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +6.22%  '
$ws.Range('D18').Value = '2.365.35'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '43.887.88'
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  +2.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.67'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '77.67'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.93%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '257.19'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.54%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.53'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.13'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +10.79%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.91%  '
$ws.Range('E28').Value = '  +13.35%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.29'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.99%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '23.10'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.86%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '175.68'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.78%  '
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.135'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.29'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.97%  '
$ws.Range('E35').Value = '  +8.33%  '
$ws.Range('E36').Value = '  +5.11%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.82'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.32%  '
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.50'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('E40').Value = '  +7.83%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.202'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +19.16%  '
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.96'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  +4.90%  '
$ws.Range('E46').Value = '  +4.71%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.50'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +12.95%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '101.99'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '54.41'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +6.67%  '
